$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 08:32:46"
$wsZhCn.Range("H3").Value = "2016-03-18 08:33:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 08:32:48"
$wsDeDe.Range("H3").Value = "2016-03-18 08:33:10"
